$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the statistics block and update values:
#   Row 8: Statistic / Value        (unchanged header)
#   Row 9: Average / 13960.00       (moved up, new value)
#   Row 10: Median / 13500.00       (moved down, new value)
#   Row 11: Difference / 460.00     (stays last, new value)
#
# Values are stored as text (not numbers), so assign them as strings.

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B11").NumberFormat = "@"

$ws.Range("A9").Value = "Average"
$ws.Range("B9").Value = "13960.00"

$ws.Range("A10").Value = "Median"
$ws.Range("B10").Value = "13500.00"

$ws.Range("A11").Value = "Difference"
$ws.Range("B11").Value = "460.00"
